$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Total" header in G1 (same bold/bordered/centered style as the other headers) ---
$ws.Range("G1").Value = "Total"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

# --- Existing rows 2 and 3 gain an (empty) Total cell, left unstyled like the rest of the row ---
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)

# --- New row 4 ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "1402-05-05"
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("B4").Value = "BTCUSDT"
$ws.Range("C4").Value = "SELL"
$ws.Range("D4").Value = 0.000005
$ws.Range("E4").Value = 500000
$ws.Range("F4").Value = "LW Strategy"
$ws.Range("G4").Value = 2.5

# --- New row 5 ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "1404-01-01"
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("B5").Value = "ETHUSDT"
$ws.Range("C5").Value = "SELL"
$ws.Range("D5").Value = 0.0002
$ws.Range("E5").Value = 5000
$ws.Range("F5").Value = "Note"
$ws.Range("G5").Value = 1
